$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we touch to Text format first,
# so values like "37.76" or "0.999" are stored as literal text
# (matching the source inline strings) rather than being parsed as numbers.
$dCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($ref in $dCells) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = "43.580.32"
$ws.Range("E2").Value = "  -5.84%  "

$ws.Range("D3").Value = "2.657.49"
$ws.Range("E3").Value = "  +2.32%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "303.50"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("D6").Value = "96.94"
$ws.Range("E6").Value = "  -2.59%  "

$ws.Range("D7").Value = "0.586"
$ws.Range("E7").Value = "  -2.26%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "0.564"
$ws.Range("E9").Value = "  -2.88%  "

$ws.Range("D10").Value = "37.47"
$ws.Range("E10").Value = "  -4.59%  "

$ws.Range("D11").Value = "0.0818"
$ws.Range("E11").Value = "  -2.85%  "

$ws.Range("D12").Value = "7.89"
$ws.Range("E12").Value = "  -3.75%  "

$ws.Range("D13").Value = "3.064.40"
$ws.Range("E13").Value = "  +2.43%  "

$ws.Range("E14").Value = "  +1.56%  "

$ws.Range("D15").Value = "2.665.44"
$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("D16").Value = "0.902"
$ws.Range("E16").Value = "  -1.62%  "

$ws.Range("D17").Value = "14.74"
$ws.Range("E17").Value = "  -1.22%  "

$ws.Range("D18").Value = "43.640.89"
$ws.Range("E18").Value = "  -6.03%  "

$ws.Range("D19").Value = "6.76"
$ws.Range("E19").Value = "  +1.61%  "

$ws.Range("D20").Value = "0.0₃0987"
$ws.Range("E20").Value = "  -2.13%  "

$ws.Range("D21").Value = "12.56"
$ws.Range("E21").Value = "  -2.88%  "

$ws.Range("D22").Value = "74.12"
$ws.Range("E22").Value = "  +3.17%  "

$ws.Range("D23").Value = "272.12"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").Value = "2.26"
$ws.Range("E24").Value = "  +4.21%  "

$ws.Range("D25").Value = "2.98"
$ws.Range("E25").Value = "  -1.29%  "

$ws.Range("D26").Value = "30.16"
$ws.Range("E26").Value = "  -0.28%  "

$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").Value = "10.35"
$ws.Range("E28").Value = "  -1.64%  "

$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  -3.19%  "

$ws.Range("D30").Value = "37.76"
$ws.Range("E30").Value = "  -3.39%  "

$ws.Range("D31").Value = "6.09"
$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").Value = "3.69"
$ws.Range("E32").Value = "  +1.77%  "

$ws.Range("D33").Value = "2.30"
$ws.Range("E33").Value = "  +5.67%  "

$ws.Range("D34").Value = "153.33"
$ws.Range("E34").Value = "  +2.19%  "

$ws.Range("E35").Value = "  -1.99%  "

$ws.Range("D36").Value = "0.0824"
$ws.Range("E36").Value = "  -1.76%  "

$ws.Range("E37").Value = "  -2.56%  "

$ws.Range("D38").Value = "25.31"
$ws.Range("E38").Value = "  +9.40%  "

$ws.Range("E39").Value = "  -0.47%  "

$ws.Range("D40").Value = "15.86"
$ws.Range("E40").Value = "  -1.37%  "

$ws.Range("D41").Value = "3.54"
$ws.Range("E41").Value = "  -1.49%  "

$ws.Range("D42").Value = "0.0319"
$ws.Range("E42").Value = "  -2.84%  "

$ws.Range("D43").Value = "3.89"
$ws.Range("E43").Value = "  -5.07%  "

$ws.Range("D44").Value = "2.117.49"
$ws.Range("E44").Value = "  -0.81%  "

$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").Value = "90.19"
$ws.Range("E46").Value = "  -3.70%  "

$ws.Range("D47").Value = "9.21"
$ws.Range("E47").Value = "  -4.83%  "

$ws.Range("D48").Value = "2.917.83"
$ws.Range("E48").Value = "  +2.72%  "

$ws.Range("D49").Value = "108.51"
$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("E50").Value = "  +3.75%  "

$ws.Range("D51").Value = "0.194"
$ws.Range("E51").Value = "  -2.83%  "

# Remove the temporary Text number format again so the cells end up
# with no explicit style, exactly like the rest of the sheet.
foreach ($ref in $dCells) { $ws.Range($ref).ClearFormats() }